$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VM")
$ws.Columns("J:J").Insert()
Write-Output "done"
